# posts.xlsx update: remove the "「私はサッカーが好きです」" entry (old row 194).
# All rows below it (195..253) shift up by one to close the gap, and the
# sheet's used range shrinks from A1:C253 to A1:C252.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(194).Delete()
